$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("results")
$ws2 = $wb.Worksheets.Item("stats")

# ---------------------------------------------------------------------------
# Sheet "results": insert a new "S*-unmerged" column before the old "var"
# column (shifting the later headers right), and append a new trailing
# "S*-unmergedND" column at the end. Mirror the same layout change in the
# row-2 data, also appending a second new trailing boolean column.
# ---------------------------------------------------------------------------
$ws1.Columns.Item(7).Insert()
$ws1.Range("G1").Value2 = "S*-unmerged"
$ws1.Range("G1").Font.Bold = $true
$ws1.Range("G1").HorizontalAlignment = -4108
$ws1.Range("G1").VerticalAlignment = -4160
$ws1.Range("G1").Borders.LineStyle = 1

$ws1.Range("M1").Value2 = "S*-unmergedND"
$ws1.Range("M1").Font.Bold = $true
$ws1.Range("M1").HorizontalAlignment = -4108
$ws1.Range("M1").VerticalAlignment = -4160
$ws1.Range("M1").Borders.LineStyle = 1

$ws1.Range("G2").Value2 = 141
$ws1.Range("H2").Value2 = 0
$ws1.Range("I2").Value2 = $false
$ws1.Range("J2").Value2 = $false
$ws1.Range("K2").Value2 = $false
$ws1.Range("L2").Value2 = $false
$ws1.Range("M2").Value2 = $false

# ---------------------------------------------------------------------------
# Sheet "stats": insert a new row for "S*-unmerged" into each of the two
# blocks (the "run 0" block and the "Average" block), right before the
# trailing "Kruskal" summary row, then refresh all of the numeric data to
# the new values.
# ---------------------------------------------------------------------------
$ws2.Rows.Item(6).Insert()
$ws2.Range("A6:B6").Font.Bold = $true
$ws2.Range("A6:B6").HorizontalAlignment = -4108
$ws2.Range("A6:B6").VerticalAlignment = -4160
$ws2.Range("A6:B6").Borders.LineStyle = 1

$ws2.Rows.Item(12).Insert()
$ws2.Range("A12:B12").Font.Bold = $true
$ws2.Range("A12:B12").HorizontalAlignment = -4108
$ws2.Range("A12:B12").VerticalAlignment = -4160
$ws2.Range("A12:B12").Borders.LineStyle = 1

$ws2.Range("A2").Value2 = "run 0"
$ws2.Range("B2").Value2 = "S*-BS"
$ws2.Range("C2").Value2 = 58
$ws2.Range("D2").Value2 = 0.0001062355004251003
$ws2.Range("E2").Value2 = 0.02906543994322419
$ws2.Range("F2").Value2 = 58
$ws2.Range("G2").Value2 = 0.003428105264902115
$ws2.Range("H2").Value2 = 0.005221989937126637
$ws2.Range("I2").Value2 = 0.004946233239024878
$ws2.Range("J2").Value2 = 0.01074812281876802
$ws2.Range("K2").Value2 = 0.00130230700597167

$ws2.Range("B3").Value2 = "S*-HS"
$ws2.Range("C3").Value2 = 58
$ws2.Range("D3").Value2 = 0.001255751121789217
$ws2.Range("E3").Value2 = 0.02114645065739751
$ws2.Range("F3").Value2 = 58
$ws2.Range("G3").Value2 = 0.002091515343636274
$ws2.Range("H3").Value2 = 0.004581128712743521
$ws2.Range("I3").Value2 = 0.003499050159007311
$ws2.Range("J3").Value2 = 0.008156916592270136
$ws2.Range("K3").Value2 = 0.0007926160469651222

$ws2.Range("B4").Value2 = "S*-MM"
$ws2.Range("C4").Value2 = 58
$ws2.Range("D4").Value2 = 0.002126702573150396
$ws2.Range("E4").Value2 = 0.03207659209147096
$ws2.Range("F4").Value2 = 58
$ws2.Range("G4").Value2 = 0.003157577011734247
$ws2.Range("H4").Value2 = 0.007013775408267975
$ws2.Range("I4").Value2 = 0.005517114885151386
$ws2.Range("J4").Value2 = 0.01198764285072684
$ws2.Range("K4").Value2 = 0.001156458631157875

$ws2.Range("B5").Value2 = "S*-MM0"
$ws2.Range("C5").Value2 = 58
$ws2.Range("D5").Value2 = 0.0001610321924090385
$ws2.Range("E5").Value2 = 0.02473008818924427
$ws2.Range("F5").Value2 = 58
$ws2.Range("G5").Value2 = 0.002894910518079996
$ws2.Range("H5").Value2 = 0.004614257253706455
$ws2.Range("I5").Value2 = 0.004969678353518248
$ws2.Range("J5").Value2 = 0.008209006860852242
$ws2.Range("K5").Value2 = 0.001077517401427031

$ws2.Range("B6").Value2 = "S*-unmerged"
$ws2.Range("C6").Value2 = 126
$ws2.Range("D6").Value2 = 0.002120027784258127
$ws2.Range("E6").Value2 = 0.06672624731436372
$ws2.Range("F6").Value2 = 126
$ws2.Range("G6").Value2 = 0.004351237323135138
$ws2.Range("H6").Value2 = 0.009801749140024185
$ws2.Range("I6").Value2 = 0.034101368393749
$ws2.Range("J6").Value2 = 0.01174131548032165
$ws2.Range("K6").Value2 = 0.001956654246896505

$ws2.Range("B7").Value2 = "Kruskal"
$ws2.Range("C7").Value2 = 975
$ws2.Range("E7").Value2 = 0.02039748523384333

$ws2.Range("A8").Value2 = "Average"
$ws2.Range("B8").Value2 = "S*-BS"
$ws2.Range("C8").Value2 = 58
$ws2.Range("D8").Value2 = 0.0001062355004251003
$ws2.Range("E8").Value2 = 0.02906543994322419
$ws2.Range("F8").Value2 = 58
$ws2.Range("G8").Value2 = 0.003428105264902115
$ws2.Range("H8").Value2 = 0.005221989937126637
$ws2.Range("I8").Value2 = 0.004946233239024878
$ws2.Range("J8").Value2 = 0.01074812281876802
$ws2.Range("K8").Value2 = 0.00130230700597167

$ws2.Range("B9").Value2 = "S*-HS"
$ws2.Range("C9").Value2 = 58
$ws2.Range("D9").Value2 = 0.001255751121789217
$ws2.Range("E9").Value2 = 0.02114645065739751
$ws2.Range("F9").Value2 = 58
$ws2.Range("G9").Value2 = 0.002091515343636274
$ws2.Range("H9").Value2 = 0.004581128712743521
$ws2.Range("I9").Value2 = 0.003499050159007311
$ws2.Range("J9").Value2 = 0.008156916592270136
$ws2.Range("K9").Value2 = 0.0007926160469651222

$ws2.Range("B10").Value2 = "S*-MM"
$ws2.Range("C10").Value2 = 58
$ws2.Range("D10").Value2 = 0.002126702573150396
$ws2.Range("E10").Value2 = 0.03207659209147096
$ws2.Range("F10").Value2 = 58
$ws2.Range("G10").Value2 = 0.003157577011734247
$ws2.Range("H10").Value2 = 0.007013775408267975
$ws2.Range("I10").Value2 = 0.005517114885151386
$ws2.Range("J10").Value2 = 0.01198764285072684
$ws2.Range("K10").Value2 = 0.001156458631157875

$ws2.Range("B11").Value2 = "S*-MM0"
$ws2.Range("C11").Value2 = 58
$ws2.Range("D11").Value2 = 0.0001610321924090385
$ws2.Range("E11").Value2 = 0.02473008818924427
$ws2.Range("F11").Value2 = 58
$ws2.Range("G11").Value2 = 0.002894910518079996
$ws2.Range("H11").Value2 = 0.004614257253706455
$ws2.Range("I11").Value2 = 0.004969678353518248
$ws2.Range("J11").Value2 = 0.008209006860852242
$ws2.Range("K11").Value2 = 0.001077517401427031

$ws2.Range("B12").Value2 = "S*-unmerged"
$ws2.Range("C12").Value2 = 126
$ws2.Range("D12").Value2 = 0.002120027784258127
$ws2.Range("E12").Value2 = 0.06672624731436372
$ws2.Range("F12").Value2 = 126
$ws2.Range("G12").Value2 = 0.004351237323135138
$ws2.Range("H12").Value2 = 0.009801749140024185
$ws2.Range("I12").Value2 = 0.034101368393749
$ws2.Range("J12").Value2 = 0.01174131548032165
$ws2.Range("K12").Value2 = 0.001956654246896505

$ws2.Range("B13").Value2 = "Kruskal"
$ws2.Range("C13").Value2 = 975
$ws2.Range("E13").Value2 = 0.02039748523384333

Write-Host "edit complete"
